$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 65 (pushes old 65.. down to 67..)
$ws.Rows("65:66").Insert()

# New row 65: Client
$ws.Range("E65").Value2 = "Client"
$ws.Range("G65").Value2 = 120
$ws.Range("H65").Value2 = 180

# New row 66: Server
$ws.Range("E66").Value2 = "Server"
$ws.Range("G66").Value2 = 100
$ws.Range("H66").Value2 = 120

# Fill in workload numbers that appear for the first time on shifted rows
$ws.Range("G68").Value2 = 8
$ws.Range("H68").Value2 = 6

$ws.Range("G69").Value2 = 12
$ws.Range("H69").Value2 = 10

$ws.Range("G70").Value2 = 6
$ws.Range("H70").Value2 = 6

$ws.Range("G72").Value2 = 12
$ws.Range("H72").Value2 = 12

$ws.Range("G73").Value2 = 4
$ws.Range("H73").Value2 = 5

$ws.Range("G74").Value2 = 6
$ws.Range("H74").Value2 = 4

$ws.Range("G75").Value2 = 1
$ws.Range("H75").Value2 = 1

# View state update
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("G64").Select()
